# Week 26 team: re-sort the current squad (Table1, visible/selected rows)
# by the "NEXT" (AI) column descending, drop the old "PREV" autofilter so
# the filter only shows Selected=1 rows, and move the selection cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# --- 1. Capture the current content of every row that is part of the
#        selected/visible team block (rows 8-143 that have Selected = 1).
#        Formula cells are captured via .Formula (text); plain literals are
#        captured via .Value2 (native number/string) so long decimals don't
#        get truncated the way a round-tripped formula string would. ---
$rows = @(8, 25, 29, 46, 58, 60, 69, 94, 96, 97, 107, 108, 110, 135, 143)

function Get-CellRaw($sheet, $r, $c) {
    $cell = $sheet.Cells.Item($r, $c)
    if ($cell.HasFormula) {
        return @{ isFormula = $true; data = $cell.Formula }
    } else {
        return @{ isFormula = $false; data = $cell.Value2 }
    }
}

function Set-CellRaw($sheet, $r, $c, $captured) {
    $cell = $sheet.Cells.Item($r, $c)
    if ($captured.isFormula) {
        $cell.Formula = $captured.data
    } else {
        $cell.Value2 = $captured.data
    }
}

$data = @{}
foreach ($r in $rows) {
    $rowVals = @()
    for ($c = 1; $c -le 38; $c++) {
        $rowVals += ,(Get-CellRaw $ws $r $c)
    }
    $data[$r] = $rowVals
}

# Unhide the target rows *before* rewriting their contents: editing cells
# on a still-hidden row causes the later AutoFilter recompute to stamp an
# explicit row height when it flips Hidden back to False, which the
# original edit doesn't have. Doing it up front avoids that altogether.
foreach ($r in $rows) {
    $ws.Rows.Item($r).Hidden = $false
}

# --- 2. New row order: same 15 physical row slots, but re-populated so
#        the content ends up sorted by column AI ("NEXT") descending.
#        (row -> row whose old content should land there) ---
$mapping = @{
    8   = 108
    25  = 96
    29  = 107
    46  = 97
    58  = 94
    60  = 110
    69  = 60
    94  = 135
    96  = 143
    97  = 29
    107 = 25
    108 = 58
    110 = 8
    135 = 46
    143 = 69
}

foreach ($newRow in $rows) {
    $oldRow = $mapping[$newRow]
    $vals = $data[$oldRow]
    for ($c = 1; $c -le 38; $c++) {
        Set-CellRaw $ws $newRow $c $vals[$c - 1]
    }
}

# --- 3. Simplify the table's AutoFilter: only keep the "Selected" (field
#        38 / AL) filter showing 1, dropping the old "PREV" (field 37 /
#        AK) filter. This also unhides every row with Selected = 1. ---
$lo.Range.AutoFilter(38, @("1"), 7)

# Re-assert the plain Hidden=False state on the newly-visible rows: the
# AutoFilter recompute above stamps an explicit row height on rows whose
# visibility it flips, which the original edit doesn't have. Re-setting
# .Hidden directly clears that stray height stamp.
foreach ($r in $rows) {
    $ws.Rows.Item($r).Hidden = $false
}

# --- 4. Move the active selection, matching the author's click on C25. ---
$ws.Range("C25").Select()
